$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing product row with new values
$ws.Range("C2").Value = "Arroz Diana x 1 Kilogramo"
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 2500
$ws.Range("F2").Value = 3500
$ws.Range("G2").Value = 45798.88830936589

# Add the new CategoryData worksheet right after ProductData
$cat = $wb.Worksheets.Add($null, $ws)
$cat.Name = "CategoryData"

$cat.Range("A1").Value = "id"
$cat.Range("B1").Value = "name"

$cat.Range("A2").Value = 0
$cat.Range("B2").Value = "Alimentos y bebidas"

$cat.Range("A3").Value = 1
$cat.Range("B3").Value = "Aseo personal"

$cat.Range("A4").Value = 2
$cat.Range("B4").Value = "Dulcería"

$cat.Range("A5").Value = 3
$cat.Range("B5").Value = "Limpieza"

$cat.Range("A6").Value = 4
$cat.Range("B6").Value = "Papelería"

$cat.Range("A7").Value = 5
$cat.Range("B7").Value = "Medicamentos"

$cat.Range("A8").Value = 6
$cat.Range("B8").Value = "Otros"

# Match header style (bold, bordered, centered) used in ProductData
$ws.Range("A1").Copy()
$cat.Range("A1:B1").PasteSpecial(-4122) # xlPasteFormats
